$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh — update Price (D) and Volume(1h) (E) columns
# per latest coinranking.com snapshot.

$ws.Range('D2').Value = '68.230.24'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '2.647.84'
$ws.Range('E3').Value = '  +0.09%  '
$ws.Range('E4').Value = '  -0.01%  '
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.61'
$ws.Range('D5').Style = $origStyle
$ws.Range('E5').Value = '  -0.31%  '
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.64'
$ws.Range('D6').Style = $origStyle
$ws.Range('E6').Value = '  +1.35%  '
$ws.Range('E7').Value = '  -0.01%  '
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.543'
$ws.Range('D8').Style = $origStyle
$ws.Range('E8').Value = '  -0.46%  '
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.141'
$ws.Range('D9').Style = $origStyle
$ws.Range('E9').Value = '  +2.98%  '
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('E11').Value = '  +0.64%  '
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.351'
$ws.Range('D12').Style = $origStyle
$ws.Range('E12').Value = '  +0.80%  '
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.02'
$ws.Range('D13').Style = $origStyle
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('E14').Value = '  +1.30%  '
$ws.Range('D15').Value = '3.128.84'
$ws.Range('E15').Value = '  +0.10%  '
$ws.Range('D16').Value = '68.149.17'
$ws.Range('D17').Value = '2.642.38'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('E18').Value = '  -0.19%  '
$origStyle = $ws.Range('D19').Style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '363.26'
$ws.Range('D19').Style = $origStyle
$ws.Range('E19').Value = '  -1.18%  '
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.34'
$ws.Range('D20').Style = $origStyle
$ws.Range('E20').Value = '  -1.24%  '
$ws.Range('E21').Value = '  +3.47%  '
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.80'
$ws.Range('D22').Style = $origStyle
$ws.Range('E22').Value = '  -1.23%  '
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.07'
$ws.Range('D23').Style = $origStyle
$ws.Range('E23').Value = '  -2.36%  '
$origStyle = $ws.Range('D24').Style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.29'
$ws.Range('D24').Style = $origStyle
$ws.Range('E24').Value = '  +2.74%  '
$ws.Range('E25').Value = '  +0.03%  '
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.75'
$ws.Range('D26').Style = $origStyle
$ws.Range('E26').Value = '  -2.55%  '
$ws.Range('D27').Value = '2.779.57'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('E28').Value = '  -0.49%  '
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = $origStyle
$ws.Range('E29').Value = '  +0.04%  '
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '558.21'
$ws.Range('D30').Style = $origStyle
$ws.Range('E30').Value = '  -2.73%  '
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.06'
$ws.Range('D31').Style = $origStyle
$ws.Range('E31').Value = '  +0.68%  '
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('D32').Style = $origStyle
$ws.Range('E32').Value = '  -0.98%  '
$ws.Range('E33').Value = '  +0.61%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -1.61%  '
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.55'
$ws.Range('D36').Style = $origStyle
$ws.Range('E36').Value = '  +0.45%  '
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '161.07'
$ws.Range('D37').Style = $origStyle
$ws.Range('E37').Value = '  +1.01%  '
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.69'
$ws.Range('D38').Style = $origStyle
$ws.Range('E38').Value = '  +2.41%  '
$ws.Range('E39').Value = '  +1.33%  '
$ws.Range('E40').Value = '  -2.81%  '
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.32'
$ws.Range('D41').Style = $origStyle
$ws.Range('E41').Value = '  -1.27%  '
$ws.Range('D42').Value = '0.0₆0336'
$ws.Range('E42').Value = '  +4.81%  '
$ws.Range('E43').Value = '  +0.34%  '
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.60'
$ws.Range('D44').Style = $origStyle
$ws.Range('E44').Value = '  -1.79%  '
$ws.Range('E45').Value = '  +0.04%  '
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '158.93'
$ws.Range('D46').Style = $origStyle
$ws.Range('E46').Value = '  +1.07%  '
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.73'
$ws.Range('D47').Style = $origStyle
$ws.Range('E47').Value = '  -0.28%  '
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.04'
$ws.Range('D48').Style = $origStyle
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  -1.43%  '
$ws.Range('E50').Value = '  +0.23%  '
$origStyle = $ws.Range('D51').Style
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.614'
$ws.Range('D51').Style = $origStyle
$ws.Range('E51').Value = '  -0.36%  '
